$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (columns D, E, G, H, L, M, N, O, P, Q, R, S, T)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.395701
$ws.Range("H2").Value = 1.187103
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.8315313333333334
$ws.Range("N2").Value = 2.494594
$ws.Range("O2").Value = 0.1310731488815592
$ws.Range("P2").Value = 0.1310731488815592
$ws.Range("Q2").Value = 0.3290377801313334
$ws.Range("R2").Value = 2.961340021182
$ws.Range("S2").Value = 0.1310731488815592
$ws.Range("T2").Value = 0.1310731488815592

# Replace row 3 entirely with new data (target cluster becomes FAPs)
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Ramp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.395701
$ws.Range("H3").Value = 1.187103
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.857187666666666
$ws.Range("N3").Value = 14.571563
$ws.Range("O3").Value = 0.7656318609505268
$ws.Range("P3").Value = 0.7656318609505269
$ws.Range("Q3").Value = 1.921994016887667
$ws.Range("R3").Value = 17.297946151989
$ws.Range("S3").Value = 0.7656318609505268
$ws.Range("T3").Value = 0.7656318609505269

# Add new row 4 (target cluster becomes sCs)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Calca"
$ws.Range("C4").Value = "Ramp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.395701
$ws.Range("H4").Value = 1.187103
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6553060000000001
$ws.Range("N4").Value = 1.965918
$ws.Range("O4").Value = 0.1032949901679139
$ws.Range("P4").Value = 0.1032949901679139
$ws.Range("Q4").Value = 0.259305239506
$ws.Range("R4").Value = 2.333747155554
$ws.Range("S4").Value = 0.1032949901679139
$ws.Range("T4").Value = 0.1032949901679139
